$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.045497666666666
$ws.Range("H2").Value = 9.136493
$ws.Range("I2").Value = 0.06184575966423571
$ws.Range("J2").Value = 0.06184575966423572
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.66643366666667
$ws.Range("N2").Value = 37.999301
$ws.Range("O2").Value = 0.1759291503241684
$ws.Range("P2").Value = 0.1759291503241684
$ws.Range("Q2").Value = 38.57559417682145
$ws.Range("R2").Value = 347.180347591393
$ws.Range("S2").Value = 0.01088047194888172
$ws.Range("T2").Value = 0.01088047194888172
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.045497666666666
$ws.Range("H3").Value = 9.136493
$ws.Range("I3").Value = 0.06184575966423571
$ws.Range("J3").Value = 0.06184575966423572
$ws.Range("O3").Value = 0.5164516272884614
$ws.Range("P3").Value = 0.5164516272884614
$ws.Range("Q3").Value = 113.2412016401461
$ws.Range("R3").Value = 1019.170814761315
$ws.Range("S3").Value = 0.03194034321948562
$ws.Range("T3").Value = 0.03194034321948562
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.045497666666666
$ws.Range("H4").Value = 9.136493
$ws.Range("I4").Value = 0.06184575966423571
$ws.Range("J4").Value = 0.06184575966423572
$ws.Range("M4").Value = 22.14777066666666
$ws.Range("N4").Value = 66.44331199999999
$ws.Range("O4").Value = 0.3076192223873702
$ws.Range("P4").Value = 0.3076192223873702
$ws.Range("Q4").Value = 67.45098388720176
$ws.Range("R4").Value = 607.0588549848159
$ws.Range("S4").Value = 0.01902494449586838
$ws.Range("T4").Value = 0.01902494449586838
$ws.Range("I5").Value = 0.6352626115862781
$ws.Range("J5").Value = 0.6352626115862781
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.66643366666667
$ws.Range("N5").Value = 37.999301
$ws.Range("O5").Value = 0.1759291503241684
$ws.Range("P5").Value = 0.1759291503241684
$ws.Range("Q5").Value = 396.2378800632823
$ws.Range("R5").Value = 3566.14092056954
$ws.Range("S5").Value = 0.1117612114890861
$ws.Range("T5").Value = 0.1117612114890861
$ws.Range("I6").Value = 0.6352626115862781
$ws.Range("J6").Value = 0.6352626115862781
$ws.Range("O6").Value = 0.5164516272884614
$ws.Range("P6").Value = 0.5164516272884614
$ws.Range("S6").Value = 0.3280824095092511
$ws.Range("T6").Value = 0.3280824095092511
$ws.Range("I7").Value = 0.6352626115862781
$ws.Range("J7").Value = 0.6352626115862781
$ws.Range("M7").Value = 22.14777066666666
$ws.Range("N7").Value = 66.44331199999999
$ws.Range("O7").Value = 0.3076192223873702
$ws.Range("P7").Value = 0.3076192223873702
$ws.Range("Q7").Value = 692.8379311836088
$ws.Range("R7").Value = 6235.54138065248
$ws.Range("S7").Value = 0.1954189905879408
$ws.Range("T7").Value = 0.1954189905879408
$ws.Range("G8").Value = 14.91542433333333
$ws.Range("H8").Value = 44.746273
$ws.Range("I8").Value = 0.3028916287494862
$ws.Range("J8").Value = 0.3028916287494862
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.66643366666667
$ws.Range("N8").Value = 37.999301
$ws.Range("O8").Value = 0.1759291503241684
$ws.Range("P8").Value = 0.1759291503241684
$ws.Range("Q8").Value = 188.9252329283526
$ws.Range("R8").Value = 1700.327096355173
$ws.Range("S8").Value = 0.05328746688620058
$ws.Range("T8").Value = 0.05328746688620058
$ws.Range("G9").Value = 14.91542433333333
$ws.Range("H9").Value = 44.746273
$ws.Range("I9").Value = 0.3028916287494862
$ws.Range("J9").Value = 0.3028916287494862
$ws.Range("O9").Value = 0.5164516272884614
$ws.Range("P9").Value = 0.5164516272884614
$ws.Range("Q9").Value = 554.6024851590239
$ws.Range("R9").Value = 4991.422366431215
$ws.Range("S9").Value = 0.1564288745597247
$ws.Range("T9").Value = 0.1564288745597247
$ws.Range("G10").Value = 14.91542433333333
$ws.Range("H10").Value = 44.746273
$ws.Range("I10").Value = 0.3028916287494862
$ws.Range("J10").Value = 0.3028916287494862
$ws.Range("M10").Value = 22.14777066666666
$ws.Range("N10").Value = 66.44331199999999
$ws.Range("O10").Value = 0.3076192223873702
$ws.Range("P10").Value = 0.3076192223873702
$ws.Range("Q10").Value = 330.3433975306862
$ws.Range("R10").Value = 2973.090577776176
$ws.Range("S10").Value = 0.09317528730356098
$ws.Range("T10").Value = 0.09317528730356098
